$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before U, shifting memory_consumed_bytes/pattern_count/patterns
# one column to the right (U->V, V->W, W->X)
$ws.Range("U1").EntireColumn.Insert()

# Header for the newly inserted column
$ws.Range("U1").Value = "csim"

# Values for the new "csim" column (rows 2-9)
$ws.Range("U2").Value = 1
$ws.Range("U3").Value = 0
$ws.Range("U4").Value = 1
$ws.Range("U5").Value = 1
$ws.Range("U6").Value = 0
$ws.Range("U7").Value = 1
$ws.Range("U8").Value = 1
$ws.Range("U9").Value = 1

# All "position" values (column B) are now -1
$ws.Range("B2:B9").Value = -1
